$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 6998.5
$ws.Cells.Item(34, 9).Value = 3747.875
$ws.Cells.Item(34, 11).Value = 3747.875
$ws.Cells.Item(34, 13).Value = -3544.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(36, 8).Value = 6998.5
$ws.Cells.Item(36, 9).Value = 3747.875
$ws.Cells.Item(36, 11).Value = 3747.875
$ws.Cells.Item(36, 13).Value = -3032.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 9).Value = 792.2
$ws.Cells.Item(98, 10).Value = 900
$ws.Cells.Item(98, 11).Value = 792.2
$ws.Cells.Item(98, 12).Value = 900
$ws.Cells.Item(98, 13).Value = 705.8
$ws.Cells.Item(98, 14).Value = -3896

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 9).Value = 792.2
$ws.Cells.Item(122, 10).Value = 900
$ws.Cells.Item(122, 11).Value = 2376.6
$ws.Cells.Item(122, 12).Value = 2700
$ws.Cells.Item(122, 13).Value = 73.39999999999964
$ws.Cells.Item(122, 14).Value = -7600

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 15947.625
$ws.Cells.Item(125, 9).Value = 21022
$ws.Cells.Item(125, 10).Value = 724.5
$ws.Cells.Item(125, 11).Value = 189198
$ws.Cells.Item(125, 12).Value = 6520.5
$ws.Cells.Item(125, 13).Value = -186738
$ws.Cells.Item(125, 14).Value = -11440.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5434.8623
$ws.Cells.Item(137, 9).Value = 1533.6666
$ws.Cells.Item(137, 10).Value = 6452.5654
$ws.Cells.Item(137, 11).Value = 4600.9998
$ws.Cells.Item(137, 12).Value = 19357.6962
$ws.Cells.Item(137, 13).Value = -2050.9998
$ws.Cells.Item(137, 14).Value = -24457.6962

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 16951436
$ws.Cells.Item(61, 9).Value = 22729426
$ws.Cells.Item(61, 10).Value = 2669.7334
$ws.Cells.Item(61, 11).Value = 22729426
$ws.Cells.Item(61, 12).Value = 2669.7334
$ws.Cells.Item(61, 13).Value = -22729214
$ws.Cells.Item(61, 14).Value = -3093.7334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value = 28801
$ws.Cells.Item(101, 10).Value = 28801
$ws.Cells.Item(101, 12).Value = 28801
$ws.Cells.Item(101, 14).Value = -35291

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 9585.786
$ws.Cells.Item(132, 9).Value = 6654.364
$ws.Cells.Item(132, 10).Value = 20334.334
$ws.Cells.Item(132, 11).Value = 19963.092
$ws.Cells.Item(132, 12).Value = 61003.00199999999
$ws.Cells.Item(132, 13).Value = -17433.092
$ws.Cells.Item(132, 14).Value = -66063.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 16951436
$ws.Cells.Item(136, 9).Value = 22729426
$ws.Cells.Item(136, 10).Value = 2669.7334
$ws.Cells.Item(136, 11).Value = 68188278
$ws.Cells.Item(136, 12).Value = 8009.2002
$ws.Cells.Item(136, 13).Value = -68185728
$ws.Cells.Item(136, 14).Value = -13109.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1678.3846
$ws.Cells.Item(134, 9).Value = 1632.1428
$ws.Cells.Item(134, 10).Value = 1732.3334
$ws.Cells.Item(134, 11).Value = 4896.428400000001
$ws.Cells.Item(134, 12).Value = 5197.0002
$ws.Cells.Item(134, 13).Value = -2361.428400000001
$ws.Cells.Item(134, 14).Value = -10267.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 24333.334
$ws.Cells.Item(28, 10).Value = 24333.334
$ws.Cells.Item(28, 12).Value = 24333.334
$ws.Cells.Item(28, 14).Value = -24823.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 31286990
$ws.Cells.Item(31, 9).Value = 90910860
$ws.Cells.Item(31, 10).Value = 55436
$ws.Cells.Item(31, 11).Value = 90910860
$ws.Cells.Item(31, 12).Value = 55436
$ws.Cells.Item(31, 13).Value = -90910565
$ws.Cells.Item(31, 14).Value = -56026

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 31286990
$ws.Cells.Item(34, 9).Value = 90910860
$ws.Cells.Item(34, 10).Value = 55436
$ws.Cells.Item(34, 11).Value = 90910860
$ws.Cells.Item(34, 12).Value = 55436
$ws.Cells.Item(34, 13).Value = -90910658
$ws.Cells.Item(34, 14).Value = -55840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 15000
$ws.Cells.Item(80, 10).Value = 15000
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 14).Value = -17246

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(83, 8).Value = 15000
$ws.Cells.Item(83, 10).Value = 15000
$ws.Cells.Item(83, 12).Value = 45000
$ws.Cells.Item(83, 14).Value = -56232

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1264.6
$ws.Cells.Item(122, 9).Value = 1041
$ws.Cells.Item(122, 11).Value = 3123
$ws.Cells.Item(122, 13).Value = -673

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 6500
$ws.Cells.Item(57, 9).Value = 3000
$ws.Cells.Item(57, 10).Value = 10000
$ws.Cells.Item(57, 11).Value = 9000
$ws.Cells.Item(57, 12).Value = 30000
$ws.Cells.Item(57, 13).Value = -8441
$ws.Cells.Item(57, 14).Value = -31118

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 517.3125
$ws.Cells.Item(113, 9).Value = 500
$ws.Cells.Item(113, 10).Value = 519.7857
$ws.Cells.Item(113, 11).Value = 1500
$ws.Cells.Item(113, 12).Value = 1559.3571
$ws.Cells.Item(113, 13).Value = 670
$ws.Cells.Item(113, 14).Value = -5899.3571

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 661.82104
$ws.Cells.Item(131, 10).Value = 923.55554
$ws.Cells.Item(131, 12).Value = 2770.66662
$ws.Cells.Item(131, 14).Value = -12850.66662

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 581966.3
$ws.Cells.Item(139, 9).Value = 787525.7
$ws.Cells.Item(139, 10).Value = 6400
$ws.Cells.Item(139, 11).Value = 2362577.1
$ws.Cells.Item(139, 12).Value = 19200
$ws.Cells.Item(139, 13).Value = -2357437.1
$ws.Cells.Item(139, 14).Value = -29480

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 8672
$ws.Cells.Item(95, 10).Value = 8672
$ws.Cells.Item(95, 12).Value = 8672
$ws.Cells.Item(95, 14).Value = -14164

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1924
$ws.Cells.Item(122, 9).Value = 1842.3334
$ws.Cells.Item(122, 11).Value = 5527.0002
$ws.Cells.Item(122, 13).Value = -3077.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2233.2
$ws.Cells.Item(7, 9).Value = 1503
$ws.Cells.Item(7, 10).Value = 2720
$ws.Cells.Item(7, 11).Value = 1503
$ws.Cells.Item(7, 12).Value = 2720
$ws.Cells.Item(7, 13).Value = -1391
$ws.Cells.Item(7, 14).Value = -2944

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 1625
$ws.Cells.Item(17, 9).Value = 2250
$ws.Cells.Item(17, 11).Value = 2250
$ws.Cells.Item(17, 13).Value = -2080

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2034.3478
$ws.Cells.Item(40, 9).Value = 1431.4706
$ws.Cells.Item(40, 10).Value = 3742.5
$ws.Cells.Item(40, 11).Value = 1431.4706
$ws.Cells.Item(40, 12).Value = 3742.5
$ws.Cells.Item(40, 13).Value = -1295.4706
$ws.Cells.Item(40, 14).Value = -4014.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(101, 8).Value = 13825
$ws.Cells.Item(101, 10).Value = 13825
$ws.Cells.Item(101, 12).Value = 13825
$ws.Cells.Item(101, 14).Value = -20315

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 42150
$ws.Cells.Item(123, 10).Value = 42150
$ws.Cells.Item(123, 12).Value = 42150
$ws.Cells.Item(123, 14).Value = -51950

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 2233.2
$ws.Cells.Item(126, 9).Value = 1503
$ws.Cells.Item(126, 10).Value = 2720
$ws.Cells.Item(126, 11).Value = 4509
$ws.Cells.Item(126, 12).Value = 8160
$ws.Cells.Item(126, 13).Value = -2039
$ws.Cells.Item(126, 14).Value = -13100

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 8074.95
$ws.Cells.Item(132, 9).Value = 14006.889
$ws.Cells.Item(132, 10).Value = 3221.5454
$ws.Cells.Item(132, 11).Value = 42020.667
$ws.Cells.Item(132, 12).Value = 9664.636200000001
$ws.Cells.Item(132, 13).Value = -39490.667
$ws.Cells.Item(132, 14).Value = -14724.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2133.5
$ws.Cells.Item(96, 9).Value = 1908.5
$ws.Cells.Item(96, 10).Value = 2471
$ws.Cells.Item(96, 11).Value = 1908.5
$ws.Cells.Item(96, 12).Value = 2471
$ws.Cells.Item(96, 13).Value = -535.5
$ws.Cells.Item(96, 14).Value = -5217

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 37600
$ws.Cells.Item(103, 10).Value = 37600
$ws.Cells.Item(103, 12).Value = 37600
$ws.Cells.Item(103, 14).Value = -39944

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1017.3333
$ws.Cells.Item(126, 9).Value = 876
$ws.Cells.Item(126, 10).Value = 1300
$ws.Cells.Item(126, 11).Value = 2628
$ws.Cells.Item(126, 12).Value = 3900
$ws.Cells.Item(126, 13).Value = -158
$ws.Cells.Item(126, 14).Value = -8840
